$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $c = $ws.Range($cellAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "246.55"
Set-TextValue "E2" "1.09%"
Set-TextValue "G2" "8"

# Row 3
Set-TextValue "D3" "29.86"
Set-TextValue "E3" "11.85%"
Set-TextValue "G3" "8"

# Row 4
Set-TextValue "D4" "5.151"
Set-TextValue "E4" "0.23%"
Set-TextValue "G4" "8"

# Row 5
Set-TextValue "D5" "0.05734"
Set-TextValue "E5" "2.13%"
Set-TextValue "G5" "8"

# Row 6
Set-TextValue "E6" "1.69%"
Set-TextValue "G6" "8"

# Row 7
Set-TextValue "D7" "0.8562"
Set-TextValue "E7" "4.47%"
Set-TextValue "G7" "8"

# Row 8
Set-TextValue "D8" "0.8737"
Set-TextValue "E8" "4.99%"
Set-TextValue "G8" "8"

# Row 9
Set-TextValue "E9" "1.82%"
Set-TextValue "G9" "8"

# Row 10
Set-TextValue "D10" "0.06976"
Set-TextValue "E10" "0.63%"
Set-TextValue "G10" "8"

# Row 11
Set-TextValue "D11" "0.02923"
Set-TextValue "E11" "0.91%"
Set-TextValue "G11" "8"

# Row 12
Set-TextValue "D12" "0.09367"
Set-TextValue "E12" "-0.21%"
Set-TextValue "G12" "8"

# Row 13
Set-TextValue "E13" "0.38%"
Set-TextValue "G13" "8"

# Row 14
Set-TextValue "D14" "0.04148"
Set-TextValue "E14" "-9.73%"
Set-TextValue "G14" "8"

# Row 15
Set-TextValue "D15" "0.0005988"
Set-TextValue "E15" "-94.00%"
Set-TextValue "G15" "8"

# Row 16
Set-TextValue "D16" "0.005968"
Set-TextValue "E16" "-3.51%"
Set-TextValue "G16" "8"

# Row 17
Set-TextValue "D17" "3.508"
Set-TextValue "E17" "-3.85%"
Set-TextValue "G17" "8"

# Row 18
Set-TextValue "D18" "3.020"
Set-TextValue "E18" "-0.17%"
Set-TextValue "G18" "8"

# Row 19
Set-TextValue "D19" "2.268"
Set-TextValue "E19" "-1.46%"
Set-TextValue "G19" "8"

# Row 20
Set-TextValue "D20" "0.3145"
Set-TextValue "E20" "1.05%"
Set-TextValue "G20" "8"

# Row 21
Set-TextValue "D21" "0.03313"
Set-TextValue "E21" "7.27%"
Set-TextValue "G21" "8"

# Row 22
Set-TextValue "D22" "0.1305"
Set-TextValue "E22" "1.07%"
Set-TextValue "G22" "8"

# Row 23
Set-TextValue "D23" "3.593"
Set-TextValue "E23" "-4.22%"
Set-TextValue "G23" "8"

# Row 24
Set-TextValue "E24" "2.68%"
Set-TextValue "G24" "8"

# Row 25
Set-TextValue "D25" "0.001210"
Set-TextValue "E25" "-1.21%"
Set-TextValue "G25" "8"

# Row 26
Set-TextValue "D26" "0.004495"
Set-TextValue "E26" "0.09%"
Set-TextValue "G26" "8"

# Row 27
Set-TextValue "E27" "22.57%"
Set-TextValue "G27" "8"

# Row 28
Set-TextValue "D28" "0.00007249"
Set-TextValue "E28" "-48.20%"
Set-TextValue "G28" "8"

# Row 29
Set-TextValue "G29" "8"

# Row 30
Set-TextValue "G30" "8"

# Row 31
Set-TextValue "G31" "8"

# Row 32
Set-TextValue "G32" "8"

# Row 33
Set-TextValue "G33" "8"

# Row 34
Set-TextValue "G34" "8"

# Row 35
Set-TextValue "G35" "8"

# Row 36
Set-TextValue "G36" "8"

# Row 37
Set-TextValue "G37" "8"

# Row 38
Set-TextValue "G38" "8"

# Row 39
Set-TextValue "G39" "8"

# Row 40
Set-TextValue "D40" "0.03786"
Set-TextValue "E40" "4.11%"
Set-TextValue "G40" "8"

# Row 41
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.005665"
Set-TextValue "E41" "-8.20%"
Set-TextValue "G41" "8"

# Row 42
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1067"
Set-TextValue "E42" "1.55%"
Set-TextValue "G42" "8"

# Row 43
Set-TextValue "D43" "0.002194"
Set-TextValue "E43" "-12.25%"
Set-TextValue "G43" "8"

# Row 44
Set-TextValue "D44" "0.01006"
Set-TextValue "E44" "23.98%"
Set-TextValue "G44" "8"

# Row 45
Set-TextValue "D45" "0.00005066"
Set-TextValue "E45" "-5.44%"
Set-TextValue "G45" "8"

# Row 46
Set-TextValue "E46" "-0.28%"
Set-TextValue "G46" "8"

# Row 47
Set-TextValue "D47" "0.07985"
Set-TextValue "E47" "-26.75%"
Set-TextValue "G47" "8"

# Row 48
Set-TextValue "D48" "0.002726"
Set-TextValue "E48" "6.91%"
Set-TextValue "G48" "8"

# Row 49
Set-TextValue "D49" "0.00002094"
Set-TextValue "E49" "-0.28%"
Set-TextValue "G49" "8"

# Row 50
Set-TextValue "E50" "-0.28%"
Set-TextValue "G50" "8"

# Row 51
Set-TextValue "G51" "8"
